$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.260.32'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '2.825.32'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '356.50'
$ws.Range('E5').Value = '  +2.87%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.35'
$ws.Range('E6').Value = '  -3.57%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').Value = '  +3.26%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.599'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.01'
$ws.Range('E10').Value = '  -4.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.95'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').Value = '3.267.30'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = '2.837.18'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.931'
$ws.Range('E17').Value = '  +4.16%  '
$ws.Range('D18').Value = '52.063.90'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.51'
$ws.Range('E19').Value = '  +5.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.22'
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.45'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.74'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '272.15'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  +2.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '27.02'
$ws.Range('E26').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.35'
$ws.Range('E28').Value = '  +1.33%  '
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.144'
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0486'
$ws.Range('E31').Value = '  +15.60%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '52.64'
$ws.Range('E32').Value = '  +4.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '35.10'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.96'
$ws.Range('E34').Value = '  +4.42%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.60'
$ws.Range('E35').Value = '  +12.23%  '
$ws.Range('E36').Value = '  +3.60%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +2.01%  '
$ws.Range('E39').Value = '  -3.59%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.35'
$ws.Range('E40').Value = '  -2.99%  '
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '127.70'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.54'
$ws.Range('E43').Value = '  -5.96%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.32'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.26'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').Value = '2.088.17'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('E48').Value = '  -4.16%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.94'
$ws.Range('E49').Value = '  +7.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.979'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.17'
$ws.Range('E51').Value = '  +2.23%  '
